$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wrapper element used to build the MODS update payload:
#   <update type="MODS"> ... </update>  ->  <datastream type="md_descriptive" operation="update"> ... </datastream>
$ws.Range("C1").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("BM1").Value = '</mods:mods></datastream></object>'

$wb.Save()
